# Auto update Excel log
# Appends new sensor-log rows to three mmWave sheets, matching the
# source device's log export (new readings captured 2026-02-01 ~21:33).

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        $value,
        [string]$status
    )

    # Column A holds an ISO-like date string ("2026-02-01"); Excel's COM
    # value-assignment auto-detects that pattern and silently coerces it
    # to a date serial number unless the cell is pre-formatted as Text.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location
    # E column is numeric on the mmWave(BR)/mmWave(HR) sheets and a plain
    # text status ("In Bed") on mmWave(InBed) - caller passes whichever is
    # appropriate and COM infers the right storage type either way.
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# --- mmWave(InBed): append rows 173-177 -----------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Add-LogRow $wsInBed 173 "2026-02-01" "21:33:21" "21:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 174 "2026-02-01" "21:33:22" "21:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 175 "2026-02-01" "21:33:24" "21:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 176 "2026-02-01" "21:33:25" "21:00" "Bedroom" "In Bed" "Occupied"
Add-LogRow $wsInBed 177 "2026-02-01" "21:33:27" "21:00" "Bedroom" "In Bed" "Occupied"

# --- mmWave(BR): append rows 166-168 ---------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Add-LogRow $wsBR 166 "2026-02-01" "21:33:23" "21:00" "Bedroom" 67 "Occupied"
Add-LogRow $wsBR 167 "2026-02-01" "21:33:25" "21:00" "Bedroom" 33 "Occupied"
Add-LogRow $wsBR 168 "2026-02-01" "21:33:26" "21:00" "Bedroom" 2 "Occupied"

# --- mmWave(HR): append rows 167-169 ---------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Add-LogRow $wsHR 167 "2026-02-01" "21:33:23" "21:00" "Bedroom" 115 "Occupied"
Add-LogRow $wsHR 168 "2026-02-01" "21:33:24" "21:00" "Bedroom" 81 "Occupied"
Add-LogRow $wsHR 169 "2026-02-01" "21:33:26" "21:00" "Bedroom" 50 "Occupied"

Write-Host "Appended 5 rows to mmWave(InBed), 3 rows to mmWave(BR), 3 rows to mmWave(HR)"
